$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row of data (5th row) to the register, splitting the
# "certificates" entry into a separate "document" entry (spitalizare / hospitalization).
$ws.Range("A5").Value = 4

# "11/01/2024" looks like a valid MM/DD/YYYY date, so assigning it directly
# to .Value would auto-convert it into a date serial number. Enter it as a
# formula that evaluates to the literal text, then paste back as a value so
# the cell ends up holding plain text (matching the other text dates in the
# register, e.g. B2/B3/B4), without touching the cell's number format/style.
$ws.Range("B5").Formula = '="11/01/2024"'
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial(-4163)

$ws.Range("C5").Value = "Marian"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "Calculatoare"
$ws.Range("F5").Value = "spitalizare"
